$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "37.530.45"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.73%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.036.83"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "258.29"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +5.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +0.00%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "57.56"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.37%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  -1.51%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.84"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.337.10"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.98%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.824"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.46%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.47"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.74%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.65%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.039.37"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "37.483.56"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.89%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "70.12"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0857"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.23"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.09%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "229.72"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +5.92%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -0.60%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.16"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.41%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "163.93"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.53%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.137"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.71%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "20.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -0.76%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0670"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.84%  "
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E35").Value = "  +9.30%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.44"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.81"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.24%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.39"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  +1.22%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "16.25"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.407.53"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.66%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "91.26"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  +1.47%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.39"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.62%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.26%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.227.17"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "
